$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "ID_0002"
$ws.Range("B3").Value = "André Automatizador"
$ws.Range("C3").Value = "automacaoteste"

$ws.Range("A3").Select()
